$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I5").Value = -0.7382099015725235
$ws.Range("J5").Value = 0.4519257258104675
$ws.Range("K5").Value = 0.1940638174351965
$ws.Range("L5").Value = 2.588047259604172
